$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.127.19"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "2.659.75"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("D5").Value = "'532.82"
$ws.Range("E5").Value = "  +4.13%  "
$ws.Range("D6").Value = "'156.54"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("E10").Value = "  +4.87%  "
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "3.122.43"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "61.108.09"
$ws.Range("D15").Value = "'22.10"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("D17").Value = "2.667.92"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'354.80"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "'10.71"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "'6.29"
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "'61.58"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "0.0₃0862"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("D28").Value = "'7.42"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'6.21"
$ws.Range("E30").Value = "  +6.88%  "
$ws.Range("E31").Value = "  +4.08%  "
$ws.Range("D32").Value = "'19.60"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").Value = "'150.42"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "'4.14"
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "'0.924"
$ws.Range("D37").Value = "'0.888"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'306.80"
$ws.Range("E39").Value = "  +4.69%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'3.82"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("D41").Value = "'0.652"
$ws.Range("E41").Value = "  +4.27%  "
$ws.Range("D42").Value = "'0.103"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("D44").Value = "'20.23"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'4.96"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").Value = "'19.22"
$ws.Range("E48").Value = "  +8.78%  "
$ws.Range("D49").Value = "'10.36"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").Value = "2.000.58"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "'1.85"
$ws.Range("E51").Value = "  +2.54%  "
